# Apply updated cryptocurrency price/volume data to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value to a cell, forcing text interpretation (via a leading
# apostrophe) when the text looks like a plain number, so Excel does not
# convert it to a numeric cell (which would drop formatting such as trailing
# zeros, e.g. "94.40" -> 94.4).
$numericPattern = '^\s*[-+]?(\d+\.?\d*|\.\d+)([eE][-+]?\d+)?\s*$'
function Set-TextValue($range, [string]$text) {
    if ($text -match $numericPattern) {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range('D2') '30.418.53'
Set-TextValue $ws.Range('E2') '  +0.22%  '
Set-TextValue $ws.Range('D3') '1.939.28'
Set-TextValue $ws.Range('E3') '  +0.25%  '
Set-TextValue $ws.Range('E4') '  +0.45%  '
Set-TextValue $ws.Range('D5') '0.7464'
Set-TextValue $ws.Range('E5') '  +3.39%  '
Set-TextValue $ws.Range('D6') '245.59'
Set-TextValue $ws.Range('E6') '  -2.49%  '
Set-TextValue $ws.Range('E7') '  +0.42%  '
Set-TextValue $ws.Range('D8') '27.54'
Set-TextValue $ws.Range('E8') '  -0.90%  '
Set-TextValue $ws.Range('D9') '0.3158'
Set-TextValue $ws.Range('E9') '  -4.11%  '
Set-TextValue $ws.Range('D10') '0.06952'
Set-TextValue $ws.Range('E10') '  -3.97%  '
Set-TextValue $ws.Range('D11') '0.7788'
Set-TextValue $ws.Range('E11') '  -3.31%  '
Set-TextValue $ws.Range('D12') '0.07992'
Set-TextValue $ws.Range('E12') '  -1.18%  '
Set-TextValue $ws.Range('D13') '1.935.59'
Set-TextValue $ws.Range('E13') '  +0.01%  '
Set-TextValue $ws.Range('D14') '5.348'
Set-TextValue $ws.Range('E14') '  -1.83%  '
Set-TextValue $ws.Range('D15') '94.40'
Set-TextValue $ws.Range('E15') '  -0.32%  '
Set-TextValue $ws.Range('D16') '14.39'
Set-TextValue $ws.Range('E16') '  -4.49%  '
Set-TextValue $ws.Range('D17') '30.430.69'
Set-TextValue $ws.Range('E17') '  +0.28%  '
Set-TextValue $ws.Range('D18') '252.35'
Set-TextValue $ws.Range('E18') '  -0.52%  '
Set-TextValue $ws.Range('D19') '0.000007891'
Set-TextValue $ws.Range('E19') '  -4.10%  '
Set-TextValue $ws.Range('D20') '5.743'
Set-TextValue $ws.Range('E20') '  -1.36%  '
Set-TextValue $ws.Range('D21') '2.192.06'
Set-TextValue $ws.Range('E21') '  +0.18%  '
Set-TextValue $ws.Range('E22') '  +0.34%  '
Set-TextValue $ws.Range('E23') '  +0.48%  '
Set-TextValue $ws.Range('D24') '6.665'
Set-TextValue $ws.Range('E24') '  -4.01%  '
Set-TextValue $ws.Range('D25') '9.480'
Set-TextValue $ws.Range('E25') '  -2.48%  '
Set-TextValue $ws.Range('D26') '165.61'
Set-TextValue $ws.Range('E26') '  -0.15%  '
Set-TextValue $ws.Range('D27') '18.95'
Set-TextValue $ws.Range('E27') '  -1.83%  '
Set-TextValue $ws.Range('D28') '0.1323'
Set-TextValue $ws.Range('E28') '  +2.19%  '
Set-TextValue $ws.Range('D29') '2.242'
Set-TextValue $ws.Range('E29') '  -4.39%  '
Set-TextValue $ws.Range('D30') '1.365'
Set-TextValue $ws.Range('E30') '  +0.67%  '
Set-TextValue $ws.Range('D31') '1.512'
Set-TextValue $ws.Range('E31') '  -2.29%  '
Set-TextValue $ws.Range('D32') '4.347'
Set-TextValue $ws.Range('E32') '  -2.14%  '
Set-TextValue $ws.Range('D33') '4.088'
Set-TextValue $ws.Range('E33') '  -2.55%  '
Set-TextValue $ws.Range('D34') '0.05143'
Set-TextValue $ws.Range('E34') '  -1.82%  '
Set-TextValue $ws.Range('D35') '1.273'
Set-TextValue $ws.Range('E35') '  +0.89%  '
Set-TextValue $ws.Range('D36') '0.7431'
Set-TextValue $ws.Range('E36') '  -0.94%  '
Set-TextValue $ws.Range('D37') '2.784'
Set-TextValue $ws.Range('E37') '  +0.89%  '
Set-TextValue $ws.Range('D38') '0.01943'
Set-TextValue $ws.Range('E38') '  -1.34%  '
Set-TextValue $ws.Range('D39') '2.809'
Set-TextValue $ws.Range('E39') '  +0.31%  '
Set-TextValue $ws.Range('D40') '78.07'
Set-TextValue $ws.Range('E40') '  -1.40%  '
Set-TextValue $ws.Range('D41') '6.419'
Set-TextValue $ws.Range('E41') '  -0.34%  '
Set-TextValue $ws.Range('D42') '0.4452'
Set-TextValue $ws.Range('E42') '  -1.93%  '
Set-TextValue $ws.Range('D43') '1.960'
Set-TextValue $ws.Range('E43') '  -3.56%  '
Set-TextValue $ws.Range('D44') '1.006'
Set-TextValue $ws.Range('E44') '  +0.44%  '
Set-TextValue $ws.Range('D45') '0.8319'
Set-TextValue $ws.Range('E45') '  -1.24%  '
Set-TextValue $ws.Range('D46') '101.25'
Set-TextValue $ws.Range('E46') '  -0.73%  '
Set-TextValue $ws.Range('D47') '9.762'
Set-TextValue $ws.Range('E47') '  -0.90%  '
Set-TextValue $ws.Range('D48') '7.448'
Set-TextValue $ws.Range('E48') '  -0.08%  '
Set-TextValue $ws.Range('B49') 'Elrond'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Range('D49') '37.16'
Set-TextValue $ws.Range('E49') '  +0.94%  '
Set-TextValue $ws.Range('B50') 'Maker'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue $ws.Range('D50') '975.83'
Set-TextValue $ws.Range('E50') '  +9.98%  '
Set-TextValue $ws.Range('D51') '0.06027'
Set-TextValue $ws.Range('E51') '  -0.43%  '
